# Prompt 9 ("Variantes responsivas e pseudo-classes (parcial)") is being
# marked as only partially done: the body text of the prompt, its tests,
# and its acceptance criteria get struck through (tracked as "withdrawn"),
# while the leading "Prompt: " label of the first paragraph stays as-is.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.IndexOf($needle) -ge 0) {
            return $p
        }
    }
    return $null
}

# --- Paragraph: "Prompt: No core, adicione util ..." ---------------------
# Only the text from "No core" onward (i.e. everything after the literal
# "Prompt: " label) becomes strikethrough.
$pPrompt = Find-ParagraphByText $d "No core, adicione util"
$promptStart = $pPrompt.Range.Start
$promptEnd = $pPrompt.Range.End
$promptFull = $pPrompt.Range.Text
$relIdx = $promptFull.IndexOf("No core")
$strikeStart = $promptStart + $relIdx
# Exclude the trailing paragraph mark from the strike range.
$strikeEnd = $promptEnd - 1
$rngPrompt = $d.Range($strikeStart, $strikeEnd)
$rngPrompt.Font.StrikeThrough = 1

# --- Paragraph: "Testes: - Dado decls com variants ..." ------------------
# Every run of text is struck through, but the paragraph mark itself is left
# untouched.
$pTestes = Find-ParagraphByText $d "Testes: - Dado decls com variants"
$rngTestes = $d.Range($pTestes.Range.Start, $pTestes.Range.End - 1)
$rngTestes.Font.StrikeThrough = 1

# --- Paragraph: "Critérios de aceite: - Variants funcionam ..." ----------
# Here the paragraph mark is struck through too, so use the full range.
$pCriterios = Find-ParagraphByText $d "Variants funcionam quando fornecidas"
$rngCriterios = $pCriterios.Range
$rngCriterios.Font.StrikeThrough = 1
